$p = $ppt.ActivePresentation

# --- Slide 1: title-slide subtitle "Lukas Schnüriger & Valentin Bürgler" ---
# Fix spelling of first name: "Lukas " -> "Lucas "
$s1 = $p.Slides.Item(1)
$shp1 = $s1.Shapes.Item(2)
$tr1 = $shp1.TextFrame.TextRange
$name1 = $tr1.Characters(1, 6)
$name1.Text = "Lucas "

# --- Slide 7: title "Persönliches Fazit Lukas Schnüriger" ---
# Fix spelling of first name inside the title: "Lukas " -> "Lucas "
$s7 = $p.Slides.Item(7)
$shp7 = $s7.Shapes.Item(1)
$tr7 = $shp7.TextFrame.TextRange

$part1 = $tr7.Characters(1, 13)
$part1.Text = "Persönliches "

$part2 = $tr7.Characters(14, 6)
$part2.Text = "Fazit "

$part3 = $tr7.Characters(20, 2)
$part3.Text = "Lu"

$part4 = $tr7.Characters(22, 1)
$part4.Text = "c"

$part5 = $tr7.Characters(23, 3)
$part5.Text = "as "
